$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Ignace Bleukx gave feedback on Thursday (09-Nov-2022 week), adding 4 more
# hours of processing time and a note about it.
$ws2.Range("C23").Formula = "=9+8+4+7+3+0+4"
$ws2.Range("J23").Value = "processing amazing feedback from Ignace Bleukx"

# Leave the selection where the edit happened, same as the author did.
$ws2.Activate()
$ws2.Range("J24").Select()
